$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "57.152.14"
Set-TextValue $ws.Range("E2") "  +1.53%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.015.63"
Set-TextValue $ws.Range("E3") "  +1.00%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.15%  "

# Row 5
Set-TextValue $ws.Range("D5") "517.02"
Set-TextValue $ws.Range("E5") "  +4.60%  "

# Row 6
Set-TextValue $ws.Range("D6") "139.93"
Set-TextValue $ws.Range("E6") "  +5.25%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -0.07%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.437"
Set-TextValue $ws.Range("E8") "  +3.44%  "

# Row 9
Set-TextValue $ws.Range("D9") "7.58"
Set-TextValue $ws.Range("E9") "  +5.07%  "

# Row 10
Set-TextValue $ws.Range("E10") "  +6.50%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.359"
Set-TextValue $ws.Range("E11") "  +2.65%  "

# Row 12
Set-TextValue $ws.Range("E12") "  +2.26%  "

# Row 13
Set-TextValue $ws.Range("D13") "3.527.37"
Set-TextValue $ws.Range("E13") "  +0.74%  "

# Row 14
Set-TextValue $ws.Range("D14") "25.90"
Set-TextValue $ws.Range("E14") "  +4.41%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.0000159"
Set-TextValue $ws.Range("E15") "  +11.27%  "

# Row 16
Set-TextValue $ws.Range("D16") "57.082.51"
Set-TextValue $ws.Range("E16") "  +1.52%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.006.66"
Set-TextValue $ws.Range("E17") "  +0.53%  "

# Row 18
Set-TextValue $ws.Range("D18") "6.00"
Set-TextValue $ws.Range("E18") "  +3.27%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.66"
Set-TextValue $ws.Range("E19") "  +2.97%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.92"
Set-TextValue $ws.Range("E20") "  +2.85%  "

# Row 21
Set-TextValue $ws.Range("D21") "329.72"
Set-TextValue $ws.Range("E21") "  +2.69%  "

# Row 22
Set-TextValue $ws.Range("E22") "  -0.12%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.486"
Set-TextValue $ws.Range("E23") "  +5.00%  "

# Row 24
Set-TextValue $ws.Range("D24") "63.85"
Set-TextValue $ws.Range("E24") "  +5.00%  "

# Row 25
Set-TextValue $ws.Range("D25") "0.172"
Set-TextValue $ws.Range("E25") "  +5.60%  "

# Row 26
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  +1.33%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.0₃0921"
Set-TextValue $ws.Range("E27") "  +6.95%  "

# Row 28
Set-TextValue $ws.Range("D28") "6.71"
Set-TextValue $ws.Range("E28") "  +3.33%  "

# Row 29
Set-TextValue $ws.Range("D29") "7.21"
Set-TextValue $ws.Range("E29") "  +8.44%  "

# Row 30
Set-TextValue $ws.Range("B30") "PancakeSwap"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D30") "1.82"
Set-TextValue $ws.Range("E30") "  +6.36%  "

# Row 31
Set-TextValue $ws.Range("B31") "Fetch.AI"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D31") "1.23"
Set-TextValue $ws.Range("E31") "  +5.96%  "

# Row 32
Set-TextValue $ws.Range("D32") "20.74"
Set-TextValue $ws.Range("E32") "  +5.43%  "

# Row 33
Set-TextValue $ws.Range("D33") "157.64"
Set-TextValue $ws.Range("E33") "  +4.41%  "

# Row 34
Set-TextValue $ws.Range("D34") "4.63"
Set-TextValue $ws.Range("E34") "  +4.87%  "

# Row 35
Set-TextValue $ws.Range("D35") "5.76"
Set-TextValue $ws.Range("E35") "  +0.93%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.29"
Set-TextValue $ws.Range("E36") "  -2.11%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.0682"
Set-TextValue $ws.Range("E37") "  +3.57%  "

# Row 38
Set-TextValue $ws.Range("D38") "24.16"
Set-TextValue $ws.Range("E38") "  +2.66%  "

# Row 39
Set-TextValue $ws.Range("D39") "3.044.38"
Set-TextValue $ws.Range("E39") "  +0.77%  "

# Row 40
Set-TextValue $ws.Range("D40") "37.29"
Set-TextValue $ws.Range("E40") "  +1.69%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.998"
Set-TextValue $ws.Range("E41") "  -0.21%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.297.10"
Set-TextValue $ws.Range("E42") "  +6.64%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.651"
Set-TextValue $ws.Range("E43") "  +2.43%  "

# Row 44
Set-TextValue $ws.Range("D44") "3.73"
Set-TextValue $ws.Range("E44") "  +5.33%  "

# Row 45
Set-TextValue $ws.Range("E45") "  +2.20%  "

# Row 46
Set-TextValue $ws.Range("E46") "  +0.92%  "

# Row 47
Set-TextValue $ws.Range("E47") "  +8.57%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.0242"
Set-TextValue $ws.Range("E48") "  +2.51%  "

# Row 49
Set-TextValue $ws.Range("D49") "5.89"
Set-TextValue $ws.Range("E49") "  +5.86%  "

# Row 50
Set-TextValue $ws.Range("D50") "19.44"
Set-TextValue $ws.Range("E50") "  +0.99%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.0880"
Set-TextValue $ws.Range("E51") "  +4.11%  "
